# Daily-Scrum workbook update: add the 25.03.2019 daily-scrum sheet.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the selection/active state of the previously-last sheet
#    (14.03.2019) so it is no longer the selected tab and its stored
#    selection becomes the full used range instead of a single cell.
# ---------------------------------------------------------------------------
$prevLast = $wb.Worksheets.Item("14.03.2019")
$prevLast.Activate()
$prevLast.Range("A1:E9").Select()

# ---------------------------------------------------------------------------
# 2. Add the new worksheet at the end of the workbook and name it.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "25.03.2019"

# ---------------------------------------------------------------------------
# 3. Fill in the daily-scrum table contents.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Frage"
$ws.Range("B1").Value = "Wer?"
$ws.Range("C1").Value = 43549
$ws.Range("D1").Value = "TODO"

$ws.Range("A2").Value = "Was wirst  du bis zum nächsten Mal machen?"
$ws.Range("B2").Value = "Lukas"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "1) User Stories vervollständigen"

$ws.Range("B3").Value = "Josy"
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = "2) Datenbank vervollständigen"

$ws.Range("B4").Value = "Cheda"
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = "3) GUI Desiginen"

$ws.Range("B5").Value = "Luca"
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = "4) Projekt erstellen und in GitLab einchecken"

$ws.Range("A6").Value = "Irgendwelche Hürden?"
$ws.Range("B6").Value = "Lukas"
$ws.Range("B7").Value = "Josy"
$ws.Range("B8").Value = "Cheda"
$ws.Range("B9").Value = "Luca"

# ---------------------------------------------------------------------------
# 4. Formatting: reuse the same look as the other daily-scrum sheets
#    (bold header row, centered "Wer?"/date cells, centered+merged
#    question column, date number format on C1).
# ---------------------------------------------------------------------------
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108  # xlCenter

$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = 1      # xlGeneral

$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").NumberFormat = "m/d/yyyy"

$ws.Range("D1").Font.Bold = $true

$ws.Range("A2:A9").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A2:A9").VerticalAlignment = -4108    # xlCenter

$ws.Range("C2:C5").HorizontalAlignment = -4108  # xlCenter

# ---------------------------------------------------------------------------
# 5. Merge the question columns exactly like the other sheets.
# ---------------------------------------------------------------------------
$ws.Range("A2:A5").Merge()
$ws.Range("A6:A9").Merge()

# ---------------------------------------------------------------------------
# 6. Column widths (best-fit like Excel would compute on data entry).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 41.14
$ws.Columns.Item(2).ColumnWidth = 6.57
$ws.Columns.Item(4).ColumnWidth = 41.29

# ---------------------------------------------------------------------------
# 7. Make the new sheet the active tab/selected sheet, matching the
#    recorded cursor position.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("E18").Select()
